$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New training-log entries for the 2025-12-16 session (Excel serial date 46007),
# appended after the existing last row (654).
#
# Columns: A Date | B Nom du joueur | C Volume | D Intensite | E Fatigue |
#          F Douleur | G Localisation douleur | H Plaisir | I Charge (=C*D)

$newDate = 46007

$players = @(
    "Yoann Martelat",
    "Kamal Bafounta",
    "Omar Benyounes",
    "Ilyes Boughanmi",
    "Naim Ighbane",
    "Romain Thunet",
    "Yoan Zouma",
    "Ilan Ihaddadene",
    "Levy Ndoutoume",
    "Karahali Souaré",
    "Emmanuel Valey",
    "Mattheo Haon",
    "Jeremie Laurent",
    "Sofiane Belle"
)

$volume    = @(70, 70, 70, 70, 70, 70, 70, 70, 70, 70, 70, 70, 70, 70)
$intensite = @(6, 8, 6, 7, 8, 7, 7, 8, 7, 6, 7, 7, 7, 7)
$fatigue   = @(5, 7, 7, 5, 7, 6, 9, 8, 7, 6, 3, 4, 4, 4)
$douleur   = @(6, 2, 0, 0, 7, 2, 5, 0, 0, 6, 0, 0, 0, 0)
$localisation = @(
    "Genou adducteurs ",
    "Genou",
    "",
    "",
    "Genou droit",
    "Genou",
    "",
    "",
    "",
    "Cheville adducteur gauche",
    "",
    "",
    "",
    ""
)
$plaisir   = @(5, 6, 7, 0, 8, 6, 9, 5, 5, 6, 7, 3, 6, 6)

$startRow = 655
$count = $players.Count

for ($i = 0; $i -lt $count; $i++) {
    $r = $startRow + $i

    # Copy formatting for columns A:F, H:I from the last existing data row (654).
    $ws.Range("A654:F654").Copy()
    $ws.Range("A$r`:F$r").PasteSpecial(-4122)
    $ws.Range("H654:I654").Copy()
    $ws.Range("H$r`:I$r").PasteSpecial(-4122)

    # Column G formatting depends on whether there is a "localisation douleur" value:
    # row 654 (has text) supplies the style when present, row 646 (blank) otherwise.
    if ($localisation[$i] -ne "") {
        $ws.Range("G654").Copy()
    } else {
        $ws.Range("G646").Copy()
    }
    $ws.Range("G$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $newDate
    $ws.Cells.Item($r, 2).Value = $players[$i]
    $ws.Cells.Item($r, 3).Value = $volume[$i]
    $ws.Cells.Item($r, 4).Value = $intensite[$i]
    $ws.Cells.Item($r, 5).Value = $fatigue[$i]
    $ws.Cells.Item($r, 6).Value = $douleur[$i]
    if ($localisation[$i] -ne "") {
        $ws.Cells.Item($r, 7).Value = $localisation[$i]
    }
    $ws.Cells.Item($r, 8).Value = $plaisir[$i]
    $ws.Cells.Item($r, 9).Formula = "=C$r*D$r"
}

[void]$ws.Range("K661").Select()
